# Insert a new row at position 235 (this shifts existing rows 235-309 down to 236-310
# and copies the formatting of the row above into the new blank row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("235:235").Insert()

# Populate the newly inserted row 235 with the new data record.
$ws.Cells.Item(235, 1).Value = 9
$ws.Cells.Item(235, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(235, 3).Value = "Metropolitana"
$ws.Cells.Item(235, 4).Value = 44876
$ws.Cells.Item(235, 5).Value = 13
$ws.Cells.Item(235, 6).Value = 100112030
$ws.Cells.Item(235, 7).Value = "Poroto granado"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 45
$ws.Cells.Item(235, 11).Value = 46000
$ws.Cells.Item(235, 12).Value = 48000
$ws.Cells.Item(235, 13).Value = 47111
$ws.Cells.Item(235, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(235, 15).Value = "Perú"
$ws.Cells.Item(235, 16).Value = 1884
$ws.Cells.Item(235, 17).Value = 25
$ws.Cells.Item(235, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date style used by the rest of column D.
$ws.Cells.Item(235, 4).NumberFormat = $ws.Cells.Item(236, 4).NumberFormat
